$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet name (date label moved from 07-24 to 07-25)
$ws.Name = "Through 2022-07-25"

# Update the "July" label in column A row 8 to reflect new through-date
$ws.Range("A8").Value = "July (through 07-25)"

# Update July row (row 8) values
$ws.Range("B8").Value = 30
$ws.Range("C8").Value = 46
$ws.Range("D8").Value = 56
$ws.Range("E8").Value = 59
$ws.Range("F8").Value = 39
$ws.Range("G8").Value = 112
$ws.Range("H8").Value = 124
$ws.Range("I8").Value = 139

# Update Total row (row 9) values
$ws.Range("B9").Value = 155
$ws.Range("C9").Value = 294
$ws.Range("D9").Value = 446
$ws.Range("E9").Value = 412
$ws.Range("F9").Value = 290
$ws.Range("G9").Value = 584
$ws.Range("H9").Value = 884
$ws.Range("I9").Value = 945
